$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72 (shifts existing rows 72-125 down to 73-126)
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly record
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = "Vega Monumental Concepción"
$ws.Range("C72").Value = "Bíobío"
$ws.Range("D72").Value = 44873
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = 100112001
$ws.Range("G72").Value = "Berenjena"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 220
$ws.Range("K72").Value = 10000
$ws.Range("L72").Value = 11000
$ws.Range("M72").Value = 10545
$ws.Range("N72").Value = "$/caja 60 unidades"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 176
$ws.Range("Q72").Value = 60
$ws.Range("R72").Value = "Hortaliza"
